# feat: add 2022-Q1 data
#
# 1. A new worksheet "2022-Q1" is inserted right before the "总计" sheet.
#    Its layout/style matches the other quarterly sheets (e.g. 2021-Q4), so
#    we duplicate that sheet and then patch the handful of cells that hold
#    different data.
# 2. The "总计" (totals) sheet gets a new first data row for 2022-Q1 and all
#    the other rows shift down by one, with the running index in column A
#    renumbered accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: duplicate "2021-Q4" (same headers/columns/styling as the new
# quarter sheet needs) and place the copy right before "总计".
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")
$template.Copy($totalSheet)

# The copy is inserted immediately before "总计" and becomes the active sheet.
# NOTE: worksheet variables captured before the insert resolve by *position*,
# so after the new sheet shifts everything, re-fetch "总计" by name rather
# than reusing the old $totalSheet reference.
$newSheet = $wb.ActiveSheet
$newSheet.Name = "2022-Q1"
$totalSheet = $wb.Worksheets.Item("总计")

# Helper: write a numeric-looking value as *text* (matching the source
# workbook, which stores these figures as inline strings, not numbers).
function Set-TextValue($ws, $row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# Row 2: fund 090019
Set-TextValue $newSheet 2 4 "2.31"
Set-TextValue $newSheet 2 5 "93.51"
Set-TextValue $newSheet 2 6 "1.99"
Set-TextValue $newSheet 2 7 "0.0460"
$newSheet.Cells.Item(2, 8).Value = 6

# Row 3: fund 006038
Set-TextValue $newSheet 3 4 "0.92"
Set-TextValue $newSheet 3 5 "93.51"
Set-TextValue $newSheet 3 6 "1.99"
Set-TextValue $newSheet 3 7 "0.0183"
$newSheet.Cells.Item(3, 8).Value = 6

# ---------------------------------------------------------------------
# Step 2: update "总计" - insert the 2022-Q1 totals as the new first data
# row, shifting everything else down by one row and renumbering column A.
# ---------------------------------------------------------------------
$rows = @(
    @("2022-Q1", 2, 0.06),
    @("2021-Q4", 2, 0.02),
    @("2021-Q3", 2, 0.01),
    @("2021-Q2", 4, 0.27),
    @("2021-Q1", 5, 0.2),
    @("2020-Q4", 5, 0.06)
)

# Row 7 is brand-new (the sheet previously only went down to row 6), so it
# has no formatting yet - copy the bold/centered/bordered look used by the
# other column-A index cells before filling it in.
$totalSheet.Range("A6").Copy()
$totalSheet.Range("A7").PasteSpecial(-4122)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $totalSheet.Cells.Item($r, 1).Value = $i
    $totalSheet.Cells.Item($r, 2).Value = $rows[$i][0]
    $totalSheet.Cells.Item($r, 3).Value = $rows[$i][1]
    $totalSheet.Cells.Item($r, 4).Value = $rows[$i][2]
}

Write-Output "done"
